$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = 'FAPs'
$ws.Cells.Item(2,2).Value = 'Col9a2'
$ws.Cells.Item(2,3).Value = 'Mag'
$ws.Cells.Item(2,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.1272026666666667
$ws.Cells.Item(2,8).Value = 0.381608
$ws.Cells.Item(2,9).Value = 0.4660422775506914
$ws.Cells.Item(2,10).Value = 0.5537822125863816
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.110622
$ws.Cells.Item(2,14).Value = 0.331866
$ws.Cells.Item(2,15).Value = 0.0981199486258843
$ws.Cells.Item(2,16).Value = 0.1197911327431839
$ws.Cells.Item(2,17).Value = 0.014071413392
$ws.Cells.Item(2,18).Value = 0.126642720528
$ws.Cells.Item(2,19).Value = 0.04572804433076395
$ws.Cells.Item(2,20).Value = 0.06633819853874935

$ws.Cells.Item(3,1).Value = 'FAPs'
$ws.Cells.Item(3,2).Value = 'Col9a2'
$ws.Cells.Item(3,3).Value = 'Mag'
$ws.Cells.Item(3,4).Value = 'MuSCs'
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.1272026666666667
$ws.Cells.Item(3,8).Value = 0.381608
$ws.Cells.Item(3,9).Value = 0.4660422775506914
$ws.Cells.Item(3,10).Value = 0.5537822125863816
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.6118760000000001
$ws.Cells.Item(3,14).Value = 1.223752
$ws.Cells.Item(3,15).Value = 0.5427242473053424
$ws.Cells.Item(3,16).Value = 0.4417284032613671
$ws.Cells.Item(3,17).Value = 0.07783225886933334
$ws.Cells.Item(3,18).Value = 0.4669935532160001
$ws.Cells.Item(3,19).Value = 0.2529324442961665
$ws.Cells.Item(3,20).Value = 0.2446213325203293

$ws.Cells.Item(4,1).Value = 'FAPs'
$ws.Cells.Item(4,2).Value = 'Col9a2'
$ws.Cells.Item(4,3).Value = 'Mag'
$ws.Cells.Item(4,4).Value = 'Neutrophils'
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.1272026666666667
$ws.Cells.Item(4,8).Value = 0.381608
$ws.Cells.Item(4,9).Value = 0.4660422775506914
$ws.Cells.Item(4,10).Value = 0.5537822125863816
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.192067
$ws.Cells.Item(4,14).Value = 0.576201
$ws.Cells.Item(4,15).Value = 0.1703603638763331
$ws.Cells.Item(4,16).Value = 0.2079868696333922
$ws.Cells.Item(4,17).Value = 0.02443143457866667
$ws.Cells.Item(4,18).Value = 0.219882911208
$ws.Cells.Item(4,19).Value = 0.07939513198529081
$ws.Cells.Item(4,20).Value = 0.1151794288544952

$ws.Cells.Item(5,1).Value = 'FAPs'
$ws.Cells.Item(5,2).Value = 'Col9a2'
$ws.Cells.Item(5,3).Value = 'Mag'
$ws.Cells.Item(5,4).Value = 'Resolving-Mac'
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.1272026666666667
$ws.Cells.Item(5,8).Value = 0.381608
$ws.Cells.Item(5,9).Value = 0.4660422775506914
$ws.Cells.Item(5,10).Value = 0.5537822125863816
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.212851
$ws.Cells.Item(5,14).Value = 0.638553
$ws.Cells.Item(5,15).Value = 0.18879544019244
$ws.Cells.Item(5,16).Value = 0.2304935943620568
$ws.Cells.Item(5,17).Value = 0.02707521480266667
$ws.Cells.Item(5,18).Value = 0.243676933224
$ws.Cells.Item(5,19).Value = 0.0879866569384701
$ws.Cells.Item(5,20).Value = 0.1276432526728077

$ws.Cells.Item(6,1).Value = 'MuSCs'
$ws.Cells.Item(6,2).Value = 'Col9a2'
$ws.Cells.Item(6,3).Value = 'Mag'
$ws.Cells.Item(6,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.5
$ws.Cells.Item(6,7).Value = 0.129733
$ws.Cells.Item(6,8).Value = 0.259466
$ws.Cells.Item(6,9).Value = 0.475312856073383
$ws.Cells.Item(6,10).Value = 0.3765320841568785
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.110622
$ws.Cells.Item(6,14).Value = 0.331866
$ws.Cells.Item(6,15).Value = 0.0981199486258843
$ws.Cells.Item(6,16).Value = 0.1197911327431839
$ws.Cells.Item(6,17).Value = 0.014351323926
$ws.Cells.Item(6,18).Value = 0.08610794355599999
$ws.Cells.Item(6,19).Value = 0.04663767301914268
$ws.Cells.Item(6,20).Value = 0.04510520487530434

$ws.Cells.Item(7,1).Value = 'MuSCs'
$ws.Cells.Item(7,2).Value = 'Col9a2'
$ws.Cells.Item(7,3).Value = 'Mag'
$ws.Cells.Item(7,4).Value = 'MuSCs'
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.5
$ws.Cells.Item(7,7).Value = 0.129733
$ws.Cells.Item(7,8).Value = 0.259466
$ws.Cells.Item(7,9).Value = 0.475312856073383
$ws.Cells.Item(7,10).Value = 0.3765320841568785
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.6118760000000001
$ws.Cells.Item(7,14).Value = 1.223752
$ws.Cells.Item(7,15).Value = 0.5427242473053424
$ws.Cells.Item(7,16).Value = 0.4417284032613671
$ws.Cells.Item(7,17).Value = 0.079380509108
$ws.Cells.Item(7,18).Value = 0.317522036432
$ws.Cells.Item(7,19).Value = 0.2579638120469793
$ws.Cells.Item(7,20).Value = 0.1663249163112926

$ws.Cells.Item(8,1).Value = 'MuSCs'
$ws.Cells.Item(8,2).Value = 'Col9a2'
$ws.Cells.Item(8,3).Value = 'Mag'
$ws.Cells.Item(8,4).Value = 'Neutrophils'
$ws.Cells.Item(8,5).Value = 1
$ws.Cells.Item(8,6).Value = 0.5
$ws.Cells.Item(8,7).Value = 0.129733
$ws.Cells.Item(8,8).Value = 0.259466
$ws.Cells.Item(8,9).Value = 0.475312856073383
$ws.Cells.Item(8,10).Value = 0.3765320841568785
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.192067
$ws.Cells.Item(8,14).Value = 0.576201
$ws.Cells.Item(8,15).Value = 0.1703603638763331
$ws.Cells.Item(8,16).Value = 0.2079868696333922
$ws.Cells.Item(8,17).Value = 0.024917428111
$ws.Cells.Item(8,18).Value = 0.149504568666
$ws.Cells.Item(8,19).Value = 0.08097447111576067
$ws.Cells.Item(8,20).Value = 0.07831372950032614

$ws.Cells.Item(9,1).Value = 'MuSCs'
$ws.Cells.Item(9,2).Value = 'Col9a2'
$ws.Cells.Item(9,3).Value = 'Mag'
$ws.Cells.Item(9,4).Value = 'Resolving-Mac'
$ws.Cells.Item(9,5).Value = 1
$ws.Cells.Item(9,6).Value = 0.5
$ws.Cells.Item(9,7).Value = 0.129733
$ws.Cells.Item(9,8).Value = 0.259466
$ws.Cells.Item(9,9).Value = 0.475312856073383
$ws.Cells.Item(9,10).Value = 0.3765320841568785
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.212851
$ws.Cells.Item(9,14).Value = 0.638553
$ws.Cells.Item(9,15).Value = 0.18879544019244
$ws.Cells.Item(9,16).Value = 0.2304935943620568
$ws.Cells.Item(9,17).Value = 0.027613798783
$ws.Cells.Item(9,18).Value = 0.165682792698
$ws.Cells.Item(9,19).Value = 0.08973689989150024
$ws.Cells.Item(9,20).Value = 0.08678823346995539

$ws.Cells.Item(10,1).Value = 'Neutrophils'
$ws.Cells.Item(10,2).Value = 'Col9a2'
$ws.Cells.Item(10,3).Value = 'Mag'
$ws.Cells.Item(10,4).Value = 'Inflammatory-Mac'
$ws.Cells.Item(10,5).Value = 1
$ws.Cells.Item(10,6).Value = 0.3333333333333333
$ws.Cells.Item(10,7).Value = 0.01600666666666667
$ws.Cells.Item(10,8).Value = 0.04802
$ws.Cells.Item(10,9).Value = 0.05864486637592556
$ws.Cells.Item(10,10).Value = 0.06968570325674002
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.110622
$ws.Cells.Item(10,14).Value = 0.331866
$ws.Cells.Item(10,15).Value = 0.0981199486258843
$ws.Cells.Item(10,16).Value = 0.1197911327431839
$ws.Cells.Item(10,17).Value = 0.00177068948
$ws.Cells.Item(10,18).Value = 0.01593620532
$ws.Cells.Item(10,19).Value = 0.005754231275977665
$ws.Cells.Item(10,20).Value = 0.00834772932913027

$ws.Cells.Item(11,1).Value = 'Neutrophils'
$ws.Cells.Item(11,2).Value = 'Col9a2'
$ws.Cells.Item(11,3).Value = 'Mag'
$ws.Cells.Item(11,4).Value = 'MuSCs'
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 0.3333333333333333
$ws.Cells.Item(11,7).Value = 0.01600666666666667
$ws.Cells.Item(11,8).Value = 0.04802
$ws.Cells.Item(11,9).Value = 0.05864486637592556
$ws.Cells.Item(11,10).Value = 0.06968570325674002
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.6118760000000001
$ws.Cells.Item(11,14).Value = 1.223752
$ws.Cells.Item(11,15).Value = 0.5427242473053424
$ws.Cells.Item(11,16).Value = 0.4417284032613671
$ws.Cells.Item(11,17).Value = 0.009794095173333335
$ws.Cells.Item(11,18).Value = 0.05876457104000001
$ws.Cells.Item(11,19).Value = 0.03182799096219658
$ws.Cells.Item(11,20).Value = 0.03078215442974522

$ws.Cells.Item(12,1).Value = 'Neutrophils'
$ws.Cells.Item(12,2).Value = 'Col9a2'
$ws.Cells.Item(12,3).Value = 'Mag'
$ws.Cells.Item(12,4).Value = 'Neutrophils'
$ws.Cells.Item(12,5).Value = 1
$ws.Cells.Item(12,6).Value = 0.3333333333333333
$ws.Cells.Item(12,7).Value = 0.01600666666666667
$ws.Cells.Item(12,8).Value = 0.04802
$ws.Cells.Item(12,9).Value = 0.05864486637592556
$ws.Cells.Item(12,10).Value = 0.06968570325674002
$ws.Cells.Item(12,11).Value = 2
$ws.Cells.Item(12,12).Value = 0.6666666666666666
$ws.Cells.Item(12,13).Value = 0.192067
$ws.Cells.Item(12,14).Value = 0.576201
$ws.Cells.Item(12,15).Value = 0.1703603638763331
$ws.Cells.Item(12,16).Value = 0.2079868696333922
$ws.Cells.Item(12,17).Value = 0.003074352446666666
$ws.Cells.Item(12,18).Value = 0.02766917202
$ws.Cells.Item(12,19).Value = 0.009990760775281611
$ws.Cells.Item(12,20).Value = 0.01449371127857084

$ws.Cells.Item(13,1).Value = 'Neutrophils'
$ws.Cells.Item(13,2).Value = 'Col9a2'
$ws.Cells.Item(13,3).Value = 'Mag'
$ws.Cells.Item(13,4).Value = 'Resolving-Mac'
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = 0.3333333333333333
$ws.Cells.Item(13,7).Value = 0.01600666666666667
$ws.Cells.Item(13,8).Value = 0.04802
$ws.Cells.Item(13,9).Value = 0.05864486637592556
$ws.Cells.Item(13,10).Value = 0.06968570325674002
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 0.212851
$ws.Cells.Item(13,14).Value = 0.638553
$ws.Cells.Item(13,15).Value = 0.18879544019244
$ws.Cells.Item(13,16).Value = 0.2304935943620568
$ws.Cells.Item(13,17).Value = 0.003407035006666666
$ws.Cells.Item(13,18).Value = 0.03066331506
$ws.Cells.Item(13,19).Value = 0.01107188336246969
$ws.Cells.Item(13,20).Value = 0.01606210821929369
